$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking values in column D so Excel
# does not coerce them to floating point numbers (source cells are text).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '70.640.93'
$ws.Range('E2').Value = '  +2.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.818.39'
$ws.Range('E3').Value = '  +1.17%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '684.21'
$ws.Range('E5').Value = '  +9.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.95'
$ws.Range('E6').Value = '  +2.76%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.817.22'
$ws.Range('E7').Value = '  +1.16%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  +0.80%  '
$ws.Range('E10').Value = '  +1.55%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.23'
$ws.Range('E11').Value = '  +6.94%  '
$ws.Range('E12').Value = '  +0.61%  '
$ws.Range('E13').Value = '  -0.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.83'
$ws.Range('E14').Value = '  +1.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.460.91'
$ws.Range('E15').Value = '  +1.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.816.76'
$ws.Range('E16').Value = '  +1.69%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '70.724.73'
$ws.Range('E17').Value = '  +2.45%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.70'
$ws.Range('E18').Value = '  +0.48%  '
$ws.Range('E19').Value = '  +2.32%  '
$ws.Range('E20').Value = '  +0.47%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.30'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '477.55'
$ws.Range('E22').Value = '  +2.14%  '
$ws.Range('E23').Value = '  +1.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.33'
$ws.Range('E24').Value = '  +0.42%  '
$ws.Range('E25').Value = '  -1.72%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.30'
$ws.Range('E26').Value = '  +2.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.32'
$ws.Range('E27').Value = '  +3.05%  '
$ws.Range('E28').Value = '  -2.05%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.968.88'
$ws.Range('E30').Value = '  +1.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.96'
$ws.Range('E31').Value = '  +11.14%  '
$ws.Range('E32').Value = '  +2.94%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.40'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.57'
$ws.Range('E34').Value = '  +2.79%  '
$ws.Range('E35').Value = '  +3.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.14'
$ws.Range('E36').Value = '  +2.29%  '
$ws.Range('B37').Value = 'RenzoRestakedETH'
$ws.Range('C37').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.773.70'
$ws.Range('E37').Value = '  +1.39%  '
$ws.Range('B38').Value = 'Binance-PegBSC-USD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  +0.39%  '
$ws.Range('E39').Value = '  +1.31%  '
$ws.Range('E40').Value = '  +1.88%  '
$ws.Range('E41').Value = '  +2.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.965'
$ws.Range('E42').Value = '  -0.39%  '
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.12'
$ws.Range('E44').Value = '  +11.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '46.06'
$ws.Range('E46').Value = '  +6.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '159.59'
$ws.Range('E47').Value = '  +3.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '48.17'
$ws.Range('E48').Value = '  +2.97%  '
$ws.Range('B49').Value = 'FLOKI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.000297'
$ws.Range('E49').Value = '  +9.75%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.44'
$ws.Range('E50').Value = '  +6.38%  '
$ws.Range('E51').Value = '  +1.68%  '
